$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N3").Value = 1.73
$ws.Range("O3").Value = 2.08
$ws.Range("N7").Value = 1.88
$ws.Range("O7").Value = 1.98
$ws.Range("AA15").Value = 6.1
$ws.Range("AB15").Value = 15
$ws.Range("AC15").Value = 80
$ws.Range("AD15").Value = 700
$ws.Range("AE15").Value = 6.9
$ws.Range("AG15").Value = 9.25
$ws.Range("AH15").Value = 23
$ws.Range("AI15").Value = 21
$ws.Range("AJ15").Value = 35
$ws.Range("H15").Value = 3.1
$ws.Range("I15").Value = 2.3
$ws.Range("L15").Value = 1.36
$ws.Range("M15").Value = 2.65
$ws.Range("N15").Value = 2.05
$ws.Range("O15").Value = 1.6
$ws.Range("P15").Value = 1.42
$ws.Range("Q15").Value = 2.47
$ws.Range("R15").Value = 1.82
$ws.Range("S15").Value = 1.78
$ws.Range("T15").Value = 8.25
$ws.Range("U15").Value = 15
$ws.Range("V15").Value = 10.75
$ws.Range("X15").Value = 28
$ws.Range("Y15").Value = 40
$ws.Range("Z15").Value = 8
$ws.Range("AA16").Value = 5.5
$ws.Range("AB16").Value = 13.5
$ws.Range("AC16").Value = 70
$ws.Range("AD16").Value = 600
$ws.Range("AE16").Value = 8.75
$ws.Range("AF16").Value = 17.5
$ws.Range("AG16").Value = 11.5
$ws.Range("AJ16").Value = 40
$ws.Range("G16").Value = 2.3
$ws.Range("H16").Value = 2.85
$ws.Range("I16").Value = 3.4
$ws.Range("K16").Value = 6
$ws.Range("L16").Value = 1.39
$ws.Range("M16").Value = 2.75
$ws.Range("N16").Value = 2.15
$ws.Range("O16").Value = 1.62
$ws.Range("Q16").Value = 2.57
$ws.Range("R16").Value = 1.8
$ws.Range("S16").Value = 1.9
$ws.Range("T16").Value = 7
$ws.Range("U16").Value = 11
$ws.Range("V16").Value = 8.75
$ws.Range("W16").Value = 25
$ws.Range("X16").Value = 19.5
$ws.Range("Y16").Value = 30
$ws.Range("Z16").Value = 6
$ws.Range("AA17").Value = 5.5
$ws.Range("AB17").Value = 15.5
$ws.Range("AE17").Value = 6.1
$ws.Range("AG17").Value = 9.5
$ws.Range("AI17").Value = 23
$ws.Range("AJ17").Value = 40
$ws.Range("H17").Value = 2.82
$ws.Range("I17").Value = 2.32
$ws.Range("M17").Value = 2.5
$ws.Range("N17").Value = 2.37
$ws.Range("O17").Value = 1.52
$ws.Range("P17").Value = 1.52
$ws.Range("Q17").Value = 2.37
$ws.Range("T17").Value = 8.25
$ws.Range("U17").Value = 17
$ws.Range("V17").Value = 11.75
$ws.Range("X17").Value = 35
$ws.Range("Y17").Value = 45
$ws.Range("AA18").Value = 7.9
$ws.Range("AB18").Value = 14.5
$ws.Range("AE18").Value = 14.5
$ws.Range("AF18").Value = 25
$ws.Range("AG18").Value = 13.5
$ws.Range("AH18").Value = 60
$ws.Range("AI18").Value = 35
$ws.Range("AJ18").Value = 35
$ws.Range("G18").Value = 1.7
$ws.Range("H18").Value = 3.95
$ws.Range("I18").Value = 4.05
$ws.Range("L18").Value = 1.19
$ws.Range("M18").Value = 3.7
$ws.Range("O18").Value = 2.07
$ws.Range("R18").Value = 1.62
$ws.Range("U18").Value = 9
$ws.Range("W18").Value = 13.5
$ws.Range("X18").Value = 12.5
$ws.Range("AE20").Value = 7.5
$ws.Range("G20").Value = 3.3
$ws.Range("I20").Value = 2.25
$ws.Range("K20").Value = 9
$ws.Range("O20").Value = 1.7
$ws.Range("U20").Value = 15
$ws.Range("AD22").Value = 451
$ws.Range("AJ22").Value = 41
$ws.Range("J22").Value = 1.1
$ws.Range("K22").Value = 7
$ws.Range("L22").Value = 1.4
$ws.Range("M22").Value = 2.75
$ws.Range("N22").Value = 2.35
$ws.Range("O22").Value = 1.57
$ws.Range("U22").Value = 8.5
$ws.Range("AA23").Value = 6
$ws.Range("AB23").Value = 15
$ws.Range("AE23").Value = 11
$ws.Range("AF23").Value = 21
$ws.Range("AI23").Value = 34
$ws.Range("G23").Value = 2
$ws.Range("I23").Value = 4.1
$ws.Range("K23").Value = 8.5
$ws.Range("N23").Value = 2.08
$ws.Range("O23").Value = 1.73
$ws.Range("T23").Value = 7
$ws.Range("U23").Value = 9.5
$ws.Range("V23").Value = 9
$ws.Range("W23").Value = 17
$ws.Range("X23").Value = 17
$ws.Range("Z23").Value = 8.5
$ws.Range("AE24").Value = 12
$ws.Range("AF24").Value = 23
$ws.Range("AH24").Value = 51
$ws.Range("AI24").Value = 41
$ws.Range("J24").Value = 1.08
$ws.Range("K24").Value = 8
$ws.Range("O24").Value = 1.7
$ws.Range("R24").Value = 1.91
$ws.Range("S24").Value = 1.91
$ws.Range("T24").Value = 7
$ws.Range("W24").Value = 15
$ws.Range("Z24").Value = 8
$ws.Range("J25").Value = 1.05
$ws.Range("K25").Value = 11
$ws.Range("N25").Value = 1.98
$ws.Range("O25").Value = 1.88
$ws.Range("AC26").Value = 41
$ws.Range("AD26").Value = 201
$ws.Range("G26").Value = 1.85
$ws.Range("H26").Value = 3.2
$ws.Range("J26").Value = 1.06
$ws.Range("K26").Value = 10
$ws.Range("P26").Value = 1.36
$ws.Range("Q26").Value = 3
$ws.Range("R26").Value = 1.75
$ws.Range("S26").Value = 2
$ws.Range("T26").Value = 8
$ws.Range("V26").Value = 8.5
$ws.Range("Z26").Value = 10
$ws.Range("K27").Value = 9
$ws.Range("L27").Value = 1.36
$ws.Range("M27").Value = 3
$ws.Range("N27").Value = 2.15
$ws.Range("O27").Value = 1.67
$ws.Range("AA30").Value = 5.8
$ws.Range("AB30").Value = 13
$ws.Range("AC30").Value = 55
$ws.Range("AD30").Value = 400
$ws.Range("AE30").Value = 9
$ws.Range("AF30").Value = 17
$ws.Range("AG30").Value = 11
$ws.Range("AH30").Value = 45
$ws.Range("AI30").Value = 29
$ws.Range("AJ30").Value = 35
$ws.Range("G30").Value = 1.78
$ws.Range("H30").Value = 3.4
$ws.Range("I30").Value = 3.9
$ws.Range("M30").Value = 3.2
$ws.Range("N30").Value = 1.88
$ws.Range("O30").Value = 1.72
$ws.Range("R30").Value = 1.83
$ws.Range("U30").Value = 7
$ws.Range("V30").Value = 7.1
$ws.Range("W30").Value = 11.75
$ws.Range("X30").Value = 12
$ws.Range("Z30").Value = 9.5
$ws.Range("AA31").Value = 6
$ws.Range("AB31").Value = 17
$ws.Range("AC31").Value = 51
$ws.Range("AD31").Value = 351
$ws.Range("AF31").Value = 17
$ws.Range("AG31").Value = 13
$ws.Range("AH31").Value = 41
$ws.Range("AI31").Value = 34
$ws.Range("AJ31").Value = 41
$ws.Range("G31").Value = 2.05
$ws.Range("I31").Value = 3.5
$ws.Range("J31").Value = 1.08
$ws.Range("K31").Value = 7.5
$ws.Range("L31").Value = 1.4
$ws.Range("M31").Value = 2.75
$ws.Range("N31").Value = 2.25
$ws.Range("O31").Value = 1.62
$ws.Range("P31").Value = 1.5
$ws.Range("Q31").Value = 2.5
$ws.Range("R31").Value = 2
$ws.Range("S31").Value = 1.73
$ws.Range("T31").Value = 6.5
$ws.Range("U31").Value = 9
$ws.Range("W31").Value = 19
$ws.Range("Y31").Value = 34
$ws.Range("Z31").Value = 8
$ws.Range("AA36").Value = 8.5
$ws.Range("AD36").Value = 67
$ws.Range("AG36").Value = 11
$ws.Range("AH36").Value = 26
$ws.Range("G36").Value = 2.55
$ws.Range("H36").Value = 3.75
$ws.Range("I36").Value = 2.4
$ws.Range("J36").Value = 1.02
$ws.Range("K36").Value = 12
$ws.Range("L36").Value = 1.1
$ws.Range("M36").Value = 6.5
$ws.Range("N36").Value = 1.4
$ws.Range("O36").Value = 2.75
$ws.Range("R36").Value = 1.36
$ws.Range("S36").Value = 3
$ws.Range("T36").Value = 17
$ws.Range("X36").Value = 17
$ws.Range("AA38").Value = 9
$ws.Range("I38").Value = 6.6
